# Lara_Gonzalo_1.1_APT122_AutoevaluacionCompetenciasFase1.docx
# "Documentos Gonzalo Lara corregidos"
#
# The autoevaluation table (3rd table in the document) gets two small
# text clean-ups (collapsing runs that were previously split across
# multiple <w:r> elements into a single run) and two whole rows removed
# (competencies that were dropped from the self-assessment).

$d = $word.ActiveDocument
$t = $d.Tables.Item(3)

# --- 1. "He diseñado modelos ERD  y eh aplicado normalización en SQL." ---
# Originally split across three runs ("...ERD  y" / " eh " / "aplicado...").
# A Find/Replace across the whole cell range collapses it back into one run.
$cell = $t.Cell(6, 7).Range
$cell.Find.ClearFormatting()
$cell.Find.Execute(
    "He diseñado modelos ERD  y eh aplicado normalización en SQL.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "He diseñado modelos ERD  y eh aplicado normalización en SQL.", 2) | Out-Null

# --- 2. "Implementar soluciones sistémicas integrales." ---
# Originally split across two runs ("Implementar" / " soluciones..."),
# merge back into a single run the same way.
$cell2 = $t.Cell(11, 1).Range
$cell2.Find.ClearFormatting()
$cell2.Find.Execute(
    "Implementar soluciones sistémicas integrales.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Implementar soluciones sistémicas integrales.", 2) | Out-Null

# --- 3. Remove the "Resolver vulnerabilidades sistémicas..." row ---
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    try {
        $rowText = $t.Cell($i, 1).Range.Text
    } catch {
        continue
    }
    if ($rowText -like "Resolver vulnerabilidades sistémicas*") {
        $t.Rows.Item($i).Delete()
        break
    }
}

# --- 4. Remove the "Transformar grandes volúmenes..." row ---
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    try {
        $rowText = $t.Cell($i, 1).Range.Text
    } catch {
        continue
    }
    if ($rowText -like "Transformar grandes volúmenes*") {
        $t.Rows.Item($i).Delete()
        break
    }
}
